# Applies the "Cyborgs" -> "Chemistry" essay rewrite described in the commit.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

function Find-Range($text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND (range): $text"
        return $null
    }
    return $r
}

# ---- Title / byline / contact info ----
Replace-Text "Interweaving Technology and Biology: Unveiling Cyborgs" "Unraveling the Enigma of Chemistry: A Journey into the Realm of Elements"
Replace-Text " Sarah Rodriguez" " Emily Carter"
Replace-Text "sarah" "emily"
Replace-Text "rodriguez@biotech" "carter@eduworld"

# ---- First body paragraph ----
Replace-Text "In the tapestry of human ingenuity, the convergence of technology and biology weaves a new thread--the cyborg" "Stepping into the realm of chemistry is akin to embarking on an enthralling odyssey, where the intricacies of matter unfold like a captivating tapestry"

Replace-Text " These entities, part organic and part machine, blur the lines between natural and artificial, challenging our understanding of life, identity, and the boundaries of human capability" " At the heart of this enigmatic discipline lie the fundamental building blocks of the universe: elements"

Replace-Text " As technology extends its reach into our bodies and biological processes, we venture into uncharted territories, where the integration of human and machine raises profound questions and promises extraordinary possibilities" " These enigmatic entities, each possessing unique properties and characteristics, orchestrate an intricate dance, giving rise to the kaleidoscope of substances that shape our world"

# New sentence pair inserted right after the run above (before the double <br/>)
$r = Find-Range "kaleidoscope of substances that shape our world."
if ($r -ne $null) {
    $r.Collapse(0)
    $r.InsertAfter(" As we embark on this transformative journey, our understanding of the microscopic realm expands, revealing the profound interconnectedness of all matter.")
}

Replace-Text "From the implantation of cochlear implants that restore hearing to the creation of biomechanical exoskeletons that augment physical abilities, we witness the seamless merging of human physiology and technological innovation" "Unraveling the enigmatic tapestry of chemistry unveils a symphony of processes, reactions, and interactions"

Replace-Text " Fueled by advances in bioengineering, nanotechnology, and artificial intelligence, cyborgization unveils a future where humans can transcend limitations, enhance their cognitive capabilities, and push the boundaries of human performance" " Atoms, the smallest units of matter, engage in a ceaseless dance, colliding, bonding, and splitting apart, orchestrating transformations that breathe life into our world"

# Two new sentence pairs inserted after the run above (before the next double <br/>)
$r = Find-Range "orchestrating transformations that breathe life into our world."
if ($r -ne $null) {
    $r.Collapse(0)
    $r.InsertAfter(" Chemistry governs the ebb and flow of energy, the intricate ballet of reactions that fuel life's processes, and the dynamic equilibrium that underpins the stability of our surroundings. Delving into this realm, we unlock the secrets of matter, piece by intricate piece, until the enigmatic enigma reveals its captivating beauty.")
}

Replace-Text "This evolutionary leap carries with it ethical, social, and philosophical implications that demand contemplation" "The study of chemistry mirrors the exploration of a hidden realm, where invisible forces shape the visible world"

Replace-Text " As we embark on this transformative journey, we must navigate the tension between human autonomy and technological intervention, ensuring that our pursuits do not override the essence of what it means to be human" " The interplay of elements, their affinities and repulsions, weaves the fabric of reality"

Replace-Text " The integration of technology must serve humanity, empowering individuals and creating a more inclusive and equitable society" " This intricate ballet of atoms gives rise to the dazzling array of materials that surround us: from the sturdy steel that forms our infrastructure to the delicate petals of a blooming flower"

# New sentence pair inserted after the run above (end of paragraph)
$r = Find-Range "delicate petals of a blooming flower."
if ($r -ne $null) {
    $r.Collapse(0)
    $r.InsertAfter(" Chemistry illuminates the intricate dance of matter, revealing the profound elegance that underpins the universe's construction.")
}

# ---- Summary heading paragraph stays the same ----

# ---- Summary body paragraph ----
Replace-Text "The emergence of cyborgs represents a pivotal moment in human history, where the intersection of technology and biology redefines the boundaries of human potential" "This essay provides an illuminating exploration into the realm of chemistry, delving into the enigmatic tapestry of elements, the intricate symphony of reactions, and the underlying elegance that governs matter's behavior"

Replace-Text " This union brings forth remarkable advancements, yet it also challenges our perception of self and identity" " Through an engaging narrative, it unveils the fundamental principles that shape our world, highlighting the profound interconnectedness of all matter and the captivating beauty inherent in the study of chemistry"

# Remove the final two runs of the summary paragraph (sentence + period)
$r = Find-Range " As we navigate the intricate landscape of cyborgization, careful consideration is vital to ensure that technology enhances, rather than diminishes, the human experience, fostering a future where humans and machines coexist harmoniously, unlocking new frontiers of exploration and possibility."
if ($r -ne $null) {
    $r.Delete()
}

# Add a new empty paragraph at the end of the document body
$endRange = $d.Paragraphs.Add()

# ---- Fix the font name typo everywhere (TimesNewToman -> Times New Roman) ----
foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    if ($rng.End -gt $rng.Start) {
        $trimmed = $d.Range($rng.Start, $rng.End - 1)
        if ($trimmed.End -gt $trimmed.Start) {
            $trimmed.Font.Name = "Times New Roman"
        }
    }
}
